# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh values to the Balmung_Profits workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 23811206
$ws.Range("J32").Value = 27779472
$ws.Range("L32").Value = 27779472
$ws.Range("N32").Value = -27780124
$ws.Range("H137").Value = 2949348.5
$ws.Range("I137").Value = 5428.7827
$ws.Range("J137").Value = 9104817
$ws.Range("K137").Value = 16286.3481
$ws.Range("L137").Value = 27314451
$ws.Range("M137").Value = -13736.3481
$ws.Range("N137").Value = -27319551
$ws.Range("H138").Value = 4846.17
$ws.Range("I138").Value = 7034.4165
$ws.Range("J138").Value = 3035.2068
$ws.Range("K138").Value = 21103.2495
$ws.Range("L138").Value = 9105.6204
$ws.Range("M138").Value = -15963.2495
$ws.Range("N138").Value = -19385.6204

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1051.5385
$ws.Range("I2").Value = 777.1
$ws.Range("J2").Value = 1966.3334
$ws.Range("K2").Value = 777.1
$ws.Range("L2").Value = 1966.3334
$ws.Range("M2").Value = -664.1
$ws.Range("N2").Value = -2192.3334
$ws.Range("H74").Value = 466699.5
$ws.Range("I74").Value = 1755.8529
$ws.Range("J74").Value = 1454704.8
$ws.Range("K74").Value = 1755.8529
$ws.Range("L74").Value = 1454704.8
$ws.Range("M74").Value = -881.8529000000001
$ws.Range("N74").Value = -1456452.8
$ws.Range("H77").Value = 466699.5
$ws.Range("I77").Value = 1755.8529
$ws.Range("J77").Value = 1454704.8
$ws.Range("K77").Value = 8779.264500000001
$ws.Range("L77").Value = 7273524
$ws.Range("M77").Value = -4411.264500000001
$ws.Range("N77").Value = -7282260
$ws.Range("H110").Value = 4331.5
$ws.Range("I110").Value = 2212.4285
$ws.Range("J110").Value = 7298.2
$ws.Range("K110").Value = 2212.4285
$ws.Range("L110").Value = 7298.2
$ws.Range("M110").Value = -167.4285
$ws.Range("N110").Value = -11388.2
$ws.Range("H116").Value = 1051.5385
$ws.Range("I116").Value = 777.1
$ws.Range("J116").Value = 1966.3334
$ws.Range("K116").Value = 777.1
$ws.Range("L116").Value = 1966.3334
$ws.Range("M116").Value = 1516.9
$ws.Range("N116").Value = -6554.3334
$ws.Range("H122").Value = 864.4
$ws.Range("I122").Value = 550.9
$ws.Range("J122").Value = 1491.4
$ws.Range("K122").Value = 1652.7
$ws.Range("L122").Value = 4474.200000000001
$ws.Range("M122").Value = 797.3000000000002
$ws.Range("N122").Value = -9374.200000000001
$ws.Range("H132").Value = 2888.7932
$ws.Range("I132").Value = 2322.4285
$ws.Range("K132").Value = 6967.2855
$ws.Range("M132").Value = -4437.2855

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1051.5385
$ws.Range("I3").Value = 777.1
$ws.Range("J3").Value = 1966.3334
$ws.Range("K3").Value = 777.1
$ws.Range("L3").Value = 1966.3334
$ws.Range("M3").Value = -663.1
$ws.Range("N3").Value = -2194.3334
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H105").Value = 14918.4
$ws.Range("I105").Value = 12106.333
$ws.Range("J105").Value = 26166.666
$ws.Range("K105").Value = 12106.333
$ws.Range("L105").Value = 26166.666
$ws.Range("M105").Value = -10359.333
$ws.Range("N105").Value = -29660.666

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 20422088
$ws.Range("I16").Value = 28572922
$ws.Range("K16").Value = 28572922
$ws.Range("M16").Value = -28572635
$ws.Range("H31").Value = 2233.3408
$ws.Range("I31").Value = 2050.6667
$ws.Range("K31").Value = 2050.6667
$ws.Range("M31").Value = -1755.6667
$ws.Range("H34").Value = 2233.3408
$ws.Range("I34").Value = 2050.6667
$ws.Range("K34").Value = 2050.6667
$ws.Range("M34").Value = -1848.6667
$ws.Range("H58").Value = 1941
$ws.Range("I58").Value = 1087.8235
$ws.Range("K58").Value = 1087.8235
$ws.Range("M58").Value = -884.8235
$ws.Range("H113").Value = 20422088
$ws.Range("I113").Value = 28572922
$ws.Range("K113").Value = 28572922
$ws.Range("M113").Value = -28570752
$ws.Range("H122").Value = 4348
$ws.Range("I122").Value = 6150
$ws.Range("J122").Value = 2546
$ws.Range("K122").Value = 18450
$ws.Range("L122").Value = 7638
$ws.Range("M122").Value = -16000
$ws.Range("N122").Value = -12538
$ws.Range("H134").Value = 2194.6667
$ws.Range("I134").Value = 1940.8422
$ws.Range("J134").Value = 2797.5
$ws.Range("K134").Value = 5822.5266
$ws.Range("L134").Value = 8392.5
$ws.Range("M134").Value = -3287.5266
$ws.Range("N134").Value = -13462.5
$ws.Range("H136").Value = 1941
$ws.Range("I136").Value = 1087.8235
$ws.Range("K136").Value = 3263.4705
$ws.Range("M136").Value = -713.4704999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 675
$ws.Range("I5").Value = 675
$ws.Range("K5").Value = 2025
$ws.Range("M5").Value = -1913
$ws.Range("H68").Value = 2980.25
$ws.Range("J68").Value = 3498
$ws.Range("L68").Value = 10494
$ws.Range("N68").Value = -12116
$ws.Range("H71").Value = 2980.25
$ws.Range("J71").Value = 3498
$ws.Range("L71").Value = 31482
$ws.Range("N71").Value = -39594
$ws.Range("H107").Value = 27778574
$ws.Range("I107").Value = 775.7222
$ws.Range("K107").Value = 2327.1666
$ws.Range("M107").Value = -407.1666
$ws.Range("H114").Value = 6227.75
$ws.Range("J114").Value = 12999.6
$ws.Range("L114").Value = 38998.8
$ws.Range("N114").Value = -45506.8
$ws.Range("H118").Value = 25250
$ws.Range("I118").Value = 25250
$ws.Range("K118").Value = 75750
$ws.Range("M118").Value = -74507
$ws.Range("H120").Value = 27683.375
$ws.Range("I120").Value = 6530
$ws.Range("J120").Value = 34734.5
$ws.Range("K120").Value = 19590
$ws.Range("L120").Value = 104203.5
$ws.Range("M120").Value = -14752
$ws.Range("N120").Value = -113879.5
$ws.Range("H122").Value = 4446256.5
$ws.Range("I122").Value = 6061160
$ws.Range("K122").Value = 54550440
$ws.Range("M122").Value = -54547990
$ws.Range("H132").Value = 1880
$ws.Range("H135").Value = 675
$ws.Range("I135").Value = 675
$ws.Range("K135").Value = 6075
$ws.Range("M135").Value = -3540

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4903.2
$ws.Range("I70").Value = 4838.6665
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 4838.6665
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -4568.6665
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 4903.2
$ws.Range("I73").Value = 4838.6665
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 4838.6665
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -3902.6665
$ws.Range("N73").Value = -6872
$ws.Range("H123").Value = 38999.5
$ws.Range("J123").Value = 38999.5
$ws.Range("L123").Value = 38999.5
$ws.Range("N123").Value = -43899.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10723.6
$ws.Range("I7").Value = 6178.8
$ws.Range("K7").Value = 6178.8
$ws.Range("M7").Value = -6066.8
$ws.Range("H22").Value = 3176.074
$ws.Range("I22").Value = 1195.4286
$ws.Range("J22").Value = 3869.3
$ws.Range("K22").Value = 1195.4286
$ws.Range("L22").Value = 3869.3
$ws.Range("M22").Value = -900.4286
$ws.Range("N22").Value = -4459.3
$ws.Range("H27").Value = 3176.074
$ws.Range("I27").Value = 1195.4286
$ws.Range("J27").Value = 3869.3
$ws.Range("K27").Value = 1195.4286
$ws.Range("L27").Value = 3869.3
$ws.Range("M27").Value = -1088.4286
$ws.Range("N27").Value = -4083.3
$ws.Range("H56").Value = 9437.5
$ws.Range("I56").Value = 9437.5
$ws.Range("K56").Value = 9437.5
$ws.Range("M56").Value = -8746.5
$ws.Range("H82").Value = 2025.7142
$ws.Range("I82").Value = 2598
$ws.Range("J82").Value = 595
$ws.Range("K82").Value = 2598
$ws.Range("L82").Value = 595
$ws.Range("M82").Value = -2237
$ws.Range("N82").Value = -1317
$ws.Range("H85").Value = 2025.7142
$ws.Range("I85").Value = 2598
$ws.Range("J85").Value = 595
$ws.Range("K85").Value = 2598
$ws.Range("L85").Value = 595
$ws.Range("M85").Value = -1350
$ws.Range("N85").Value = -3091
$ws.Range("H94").Value = 39999
$ws.Range("J94").Value = 39999
$ws.Range("L94").Value = 39999
$ws.Range("N94").Value = -41351
$ws.Range("H126").Value = 10723.6
$ws.Range("I126").Value = 6178.8
$ws.Range("K126").Value = 18536.4
$ws.Range("M126").Value = -16066.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 72606.71000000001
$ws.Range("I81").Value = 999.8333
$ws.Range("J81").Value = 126311.875
$ws.Range("K81").Value = 1999.6666
$ws.Range("L81").Value = 252623.75
$ws.Range("M81").Value = -938.6666
$ws.Range("N81").Value = -254745.75
$ws.Range("H84").Value = 72606.71000000001
$ws.Range("I84").Value = 999.8333
$ws.Range("J84").Value = 126311.875
$ws.Range("K84").Value = 9998.333000000001
$ws.Range("L84").Value = 1263118.75
$ws.Range("M84").Value = -4694.333000000001
$ws.Range("N84").Value = -1273726.75
